# "fixed number of springs"
# The Hardware table treats the Spring line item as needing 2 springs per
# unit instead of 1, so the unit price must be doubled and the label
# updated to reflect that.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is the Spring hardware line (B4 label / D4 unit-price formula).
$ws.Range("B4").Value = "Springx2"
$ws.Range("D4").Formula = "=(E4/F4)*2"

# Leave the selection where the editor ended up after making the change.
$ws.Range("E16").Select() | Out-Null
